# Update "cantidad" (quantity) values in the Datos sheet to reflect
# the latest inventory counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

$updates = @{
    "C6"  = 37
    "C8"  = 12
    "C9"  = 5
    "C12" = 58
    "C13" = 23
    "C15" = 48
    "C20" = 16
    "C23" = 108
    "C24" = 135
    "C25" = 14
    "C26" = 10
    "C29" = 36
    "C31" = 23
    "C32" = 44
    "C34" = 34
    "C38" = 2
    "C39" = 10
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
